$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ES")

# --- Row 5 ---
$ws.Range("F5").Value = "ep"
$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = "mascarenhas"

# --- Row 7 ---
$ws.Range("F7").Value = "f"

# --- Row 8 ---
$ws.Range("D8").Value = "f"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "f"
$ws.Range("G8").Value = ""

# --- Row 9 ---
$ws.Range("F9").Value = "ep"
$ws.Range("G3").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "neima"

# --- Row 10 ---
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""

# --- Row 11 ---
$ws.Range("F11").Value = "f"
$ws.Range("G11").Value = ""

# --- Row 13 ---
$ws.Range("F13").Value = "f"
$ws.Range("G13").Value = ""

# --- Selection ---
$ws.Range("E8").Select()
